# feat: add 2022-Q1 data
#
# The previous "总计" (summary) sheet is renamed to "2022-Q1" and repopulated
# with that quarter's fund-holdings detail (same layout as the other
# quarterly sheets). A brand-new "总计" sheet is appended at the end with the
# summary table, now including a 2022-Q1 row at the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# The header style (bold / centered / bordered) already lives on B1:D1 from
# the old sheet; stamp the same style onto the newly-used E1:H1 header cells.
$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,5))
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,6))
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,7))
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,4).Copy($q1.Cells.Item(1,8))
$q1.Cells.Item(1,8).Value = "仓位排名"

# Data rows. Column A keeps the bold/centered/bordered "index" style (copied
# from the existing A2 cell so brand-new rows A7:A9 pick it up too); columns
# D/E/F/G hold numeric-looking text, so they're apostrophe-protected to stay
# text instead of being coerced to numbers (matches source layout).
$styleSrc = $q1.Cells.Item(2,1)

$rows = @(
    @(0,"002270","东吴安盈量化灵活配置混合","5.24","46.02","2.58","0.1352",8),
    @(1,"159855","银华中证影视主题ETF","0.96","97.27","9.07","0.0871",1),
    @(2,"290012","泰信行业精选灵活配置混合A","0.76","92.62","5.39","0.0410",8),
    @(3,"673040","西部利得行业主题优选灵活配置混合A","4.73","29.77","0.70","0.0331",10),
    @(4,"516620","国泰中证影视主题ETF","0.33","96.08","8.82","0.0291",1),
    @(5,"003132","德邦新回报灵活配置混合","0.55","79.88","4.12","0.0227",4),
    @(6,"673043","西部利得行业主题优选灵活配置混合C","2.67","29.77","0.70","0.0187",10),
    @(7,"002583","泰信行业精选灵活配置混合C","0.00","92.62","5.39",$null,8)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q1.Cells.Item($r,1)
    $styleSrc.Copy($aCell)
    $aCell.Value = $row[0]

    $q1.Cells.Item($r,2).Value = "'" + $row[1]
    $q1.Cells.Item($r,3).Value = $row[2]
    $q1.Cells.Item($r,4).Value = "'" + $row[3]
    $q1.Cells.Item($r,5).Value = "'" + $row[4]
    $q1.Cells.Item($r,6).Value = "'" + $row[5]

    if ($row[6] -ne $null) {
        $q1.Cells.Item($r,7).Value = "'" + $row[6]
    } else {
        $q1.Cells.Item($r,7).Value = 0
    }

    $q1.Cells.Item($r,8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" summary sheet at the end
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$total.Name = "总计"

# Pull the bold/centered/bordered "header" and "index" styles straight from
# the 2022-Q1 sheet (cross-sheet Copy carries the style index across without
# minting new style entries), then overwrite with this sheet's own values.
$headerStyleSrc = $q1.Cells.Item(1,2)
$indexStyleSrc = $q1.Cells.Item(2,1)

$headerStyleSrc.Copy($total.Cells.Item(1,2))
$total.Cells.Item(1,2).Value = "日期"
$headerStyleSrc.Copy($total.Cells.Item(1,3))
$total.Cells.Item(1,3).Value = "持有数量(只)"
$headerStyleSrc.Copy($total.Cells.Item(1,4))
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0,"2022-Q1",8,0.37),
    @(1,"2021-Q4",12,2.91),
    @(2,"2021-Q3",9,0.89),
    @(3,"2021-Q2",31,3.03),
    @(4,"2021-Q1",35,7.94),
    @(5,"2020-Q4",37,9.93)
)

$tr = 2
foreach ($row in $totalRows) {
    $aCell = $total.Cells.Item($tr,1)
    $indexStyleSrc.Copy($aCell)
    $aCell.Value = $row[0]

    $total.Cells.Item($tr,2).Value = $row[1]
    $total.Cells.Item($tr,3).Value = $row[2]
    $total.Cells.Item($tr,4).Value = $row[3]

    $tr = $tr + 1
}
